$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for 2022-Q4 and renumber the existing
#    sequence index in column A (0,1,2,... shift by one).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert()

# Pull style from row 3 (the row pushed down, which retains the original
# per-column formatting) onto the freshly inserted row 2 cells.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)
$summary.Range("B3:D3").Copy()
$summary.Range("B2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 18
$summary.Range("D2").Value = 2.28

$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5
$summary.Range("A8").Value = 6
$summary.Range("A9").Value = 7
$summary.Range("A10").Value = 8

# ---------------------------------------------------------------------------
# 2) Insert a brand-new "2022-Q4" sheet right before "2022-Q3" (i.e. right
#    after "总计"), holding the quarterly fund holdings detail.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q4 = $wb.Worksheets.Add($q3)
$q4.Name = "2022-Q4"

# Header row, styled like the other quarter sheets (bold + border + centred).
$hdr = $q4.Range("B1:H1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

$q4.Range("B1").Value = "基金代码"
$q4.Range("C1").Value = "基金名称"
$q4.Range("D1").Value = "基金规模"
$q4.Range("E1").Value = "股票总仓位"
$q4.Range("F1").Value = "仓位占比"
$q4.Range("G1").Value = "持有市值(亿元)"
$q4.Range("H1").Value = "仓位排名"

$q4rows = @(
    @(0,  "010714", "东方红远见价值混合A",               "16.55", "83.49", "2.97", "0.4915", 10),
    @(1,  "001556", "天弘中证500指数增强A",               "25.50", "94.27", "1.57", "0.4004", 10),
    @(2,  "502000", "西部利得中证500指数增强（LOF）A",     "18.20", "90.25", "1.64", "0.2985", 9),
    @(3,  "001557", "天弘中证500指数增强C",               "13.20", "94.27", "1.57", "0.2072", 10),
    @(4,  "162102", "金鹰中小盘精选混合",                 "3.48",  "76.23", "4.72", "0.1643", 4),
    @(5,  "001167", "金鹰科技创新股票",                   "3.17",  "91.02", "5.12", "0.1623", 4),
    @(6,  "210009", "金鹰核心资源混合",                   "3.14",  "91.78", "4.96", "0.1557", 7),
    @(7,  "009300", "西部利得中证500指数增强（LOF）C",     "4.99",  "90.25", "1.64", "0.0818", 9),
    @(8,  "014155", "国泰君安中证500指数增强A",           "7.70",  "92.93", "1.04", "0.0801", 8),
    @(9,  "210002", "金鹰红利价值混合A",                   "0.96",  "77.22", "5.39", "0.0517", 6),
    @(10, "014156", "国泰君安中证500指数增强C",           "4.81",  "92.93", "1.04", "0.0500", 8),
    @(11, "000458", "英大领先回报混合",                   "1.81",  "93.66", "2.57", "0.0465", 4),
    @(12, "008072", "景顺长城创业板综指增强",             "1.80",  "94.12", "1.76", "0.0317", 10),
    @(13, "016563", "金鹰红利价值混合C",                   "0.52",  "77.22", "5.39", "0.0280", 6),
    @(14, "004890", "中邮健康文娱灵活配置混合",           "0.42",  "92.60", "4.46", "0.0187", 7),
    @(15, "001270", "英大灵活配置混合A",                   "0.29",  "92.68", "2.55", "0.0074", 4),
    @(16, "001271", "英大灵活配置混合B",                   "0.28",  "92.68", "2.55", "0.0071", 4)
)

$r = 2
foreach ($row in $q4rows) {
    $q4.Range("A$r").Value = $row[0]
    $q4.Range("B$r").Value = "'" + $row[1]
    $q4.Range("C$r").Value = "'" + $row[2]
    $q4.Range("D$r").Value = "'" + $row[3]
    $q4.Range("E$r").Value = "'" + $row[4]
    $q4.Range("F$r").Value = "'" + $row[5]
    $q4.Range("G$r").Value = "'" + $row[6]
    $q4.Range("H$r").Value = $row[7]
    $r = $r + 1
}

# Row 19 - same as the others except the holding value is a true 0 (number).
$q4.Range("A19").Value = 17
$q4.Range("B19").Value = "'017537"
$q4.Range("C19").Value = "'东方红远见价值混合C"
$q4.Range("D19").Value = "'0.00"
$q4.Range("E19").Value = "'83.49"
$q4.Range("F19").Value = "'2.97"
$q4.Range("G19").Value = 0
$q4.Range("H19").Value = 10

# Style column A (index numbers) like the other sheets' leading column.
$q4.Range("A2:A19").Copy()
$q4.Range("B1").PasteSpecial(-4122)
$q3.Range("A2").Copy()
$q4.Range("A2:A19").PasteSpecial(-4122)
$r = 2
foreach ($row in $q4rows) {
    $q4.Range("A$r").Value = $row[0]
    $r = $r + 1
}
$q4.Range("A19").Value = 17

# ---------------------------------------------------------------------------
# 3) Restore the originally-selected tab (2020-Q4) as the active sheet, since
#    inserting the new sheet above switched the active tab to it.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
